# "Add files via upload" — refresh the daily figure in T2 and leave the
# cursor resting on it, matching the uploaded workbook snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount figure in T2 (was 175829)
$ws.Range("T2").Value = 204387

# Move/save the active selection to T2 (was T3)
$ws.Range("T2").Select()
